$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1) "Informe Inicial" (sheet1): add row 14
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Informe Inicial")
$ws1.Range("B5").Copy()
$ws1.Range("B14").PasteSpecial(-4122)   # xlPasteFormats (numFmt date, border)
$ws1.Range("B14").Value = 43170
$ws1.Range("C14").Value = "Informe final acabado"

# ---------------------------------------------------------------------------
# 2) Create "Segundo informe" and "Informe Final" sheets (title only for now,
#    so the shared-string table fills up in the right order), then a dummy
#    sheet in between so "Informe Final" lands on sheetId 5 (matching the
#    history of a previously-deleted sheet) while keeping sheetId 3 for
#    "Segundo informe".
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Primer informe")

$seg = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$seg.Name = "Segundo informe"

$dummy = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$dummy.Name = "DummyToDelete"

$fin = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$fin.Name = "Informe Final"

$wb.Worksheets.Item("DummyToDelete").Delete()

$seg = $wb.Worksheets.Item("Segundo informe")
$fin = $wb.Worksheets.Item("Informe Final")

# -- Segundo informe: title / header skeleton (mirrors sheet1/2 layout,
#    shifted up one row since there is no leading blank row) --------------
$seg.Columns.Item(3).ColumnWidth = 57.85546875

$ws1.Range("B2:C2").Copy()
$seg.Range("B1").PasteSpecial(-4122)
$seg.Rows.Item(1).RowHeight = 18.75
$seg.Range("B1").Value = "SEGUNDO INFORME"

$ws1.Range("B3:C3").Copy()
$seg.Range("B2").PasteSpecial(-4122)
$seg.Rows.Item(2).RowHeight = 15.75

$ws1.Range("B4:C4").Copy()
$seg.Range("B3").PasteSpecial(-4122)
$seg.Rows.Item(3).RowHeight = 15.75

$seg.Range("B1:C1").Merge()

# -- Informe Final: title / header skeleton --------------------------------
$fin.Columns.Item(3).ColumnWidth = 86.42578125

$ws1.Range("B2:C2").Copy()
$fin.Range("B1").PasteSpecial(-4122)
$fin.Rows.Item(1).RowHeight = 18.75
$fin.Range("B1").Value = "Informe Final"

$ws1.Range("B3:C3").Copy()
$fin.Range("B2").PasteSpecial(-4122)
$fin.Rows.Item(2).RowHeight = 15.75

$ws1.Range("B4:C4").Copy()
$fin.Range("B3").PasteSpecial(-4122)
$fin.Rows.Item(3).RowHeight = 15.75

$fin.Range("B1:C1").Merge()
$fin.PageSetup.PaperSize = 9
$fin.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 3) "Primer informe" (sheet2): fill in rows 10-14
# ---------------------------------------------------------------------------
function Set-DataRow($ws, $row, $date, $text) {
    $ws.Range("B5").Copy()
    $ws.Range("B$row").PasteSpecial(-4122)
    $ws.Range("B$row").Value = $date
    $ws.Range("C5").Copy()
    $ws.Range("C$row").PasteSpecial(-4122)
    $ws.Range("C$row").Value = $text
}

Set-DataRow $ws2 10 43208 "retoques en el desarrollo, cambios en el planning y retoques en dataset"
Set-DataRow $ws2 11 43208 "actualizado arbol de tareas"
Set-DataRow $ws2 12 43208 "actualizado gantt"
Set-DataRow $ws2 13 43208 "se han añadido imagenes"
Set-DataRow $ws2 14 43210 "version final, algunas correcciones"

# ---------------------------------------------------------------------------
# 4) "Segundo informe": rows 4-10
# ---------------------------------------------------------------------------
Set-DataRow $seg 4 43242 "Creación del documento"
Set-DataRow $seg 5 43243 "añadido progresso del desarrollo i cambios en la intro"
Set-DataRow $seg 6 43244 "se han añadido referencias y modificado el planning"
Set-DataRow $seg 7 43244 "ligeras modificaciones y correciones"
Set-DataRow $seg 8 43245 "se añade la seccion de user testing y se modifica la metodologia"

# rows 9-10 carry a different (newer) style: date fmt w/o border (B), and
# border-only w/o explicit fill (C)
$seg.Range("B9").Value = 43246
$seg.Range("B9").NumberFormat = "mm-dd-yy"
$seg.Range("C9").Value = "correcciones"
$ws1.Range("C5").Copy()
$seg.Range("C9").PasteSpecial(-4122)
$seg.Range("C9").Value = "correcciones"

$seg.Range("B10").Value = 43247
$seg.Range("B10").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$seg.Range("C10").PasteSpecial(-4122)
$seg.Range("C10").Value = "version final"

# ---------------------------------------------------------------------------
# 5) "Informe Final": rows 4-23
# ---------------------------------------------------------------------------
# rows 4-5 use date-format-with-border style (even C4 which is text)
$ws1.Range("B5").Copy()
$fin.Range("B4").PasteSpecial(-4122)
$fin.Range("B4").Value = 43268
$fin.Range("B4").NumberFormat = "mm-dd-yy"
$fin.Range("C4").Value = "creacion del documento"
$fin.Range("B4").Copy()
$fin.Range("C4").PasteSpecial(-4122)
$fin.Range("C4").Value = "creacion del documento"

$fin.Range("B4").Copy()
$fin.Range("B5").PasteSpecial(-4122)
$fin.Range("B5").Value = 43268
Set-DataRow $fin 5 43268 "se planea y se crea la estructura del documento"
$fin.Range("B5").Copy()
$fin.Range("B5").PasteSpecial(-4122)

Set-DataRow $fin 6 43269 "se reestructura el documento despues de la reunion"
Set-DataRow $fin 7 43269 "se añade la metodologia y los objetivos"
Set-DataRow $fin 8 43271 "se añade la introduccion y el state of the art"

$fin.Range("B9").Value = 43272
$fin.Range("B9").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C9").PasteSpecial(-4122)
$fin.Range("C9").Value = "correccion ortografica"

$fin.Range("B10").Value = 43272
$fin.Range("B10").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C10").PasteSpecial(-4122)
$fin.Range("C10").Value = "se añade la seccion de desarrollo de libelas"

$fin.Range("B11").Value = 43274
$fin.Range("B11").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C11").PasteSpecial(-4122)
$fin.Range("C11").Value = "se han añadido las explicaciones sobre los efectos de la distancia"

$fin.Range("B12").Value = 43276
$fin.Range("B12").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C12").PasteSpecial(-4122)
$fin.Range("C12").Value = "añadido user testing y revisado todo el texto"

$fin.Range("B13").Value = 43276
$fin.Range("B13").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C13").PasteSpecial(-4122)
$fin.Range("C13").Value = "se han añadido las graficas del user testing"

$fin.Range("B14").Value = 43276
$fin.Range("B14").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C14").PasteSpecial(-4122)
$fin.Range("C14").Value = "se ha modificado la explicacion del pipeline"

$fin.Range("B15").Value = 43276
$fin.Range("B15").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C15").PasteSpecial(-4122)
$fin.Range("C15").Value = "añadidas las explicaciones de la segunda session de user testing y el protocolo del user testing"

$fin.Range("B16").Value = 43277
$fin.Range("B16").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C16").PasteSpecial(-4122)
$fin.Range("C16").Value = "añadidas projectos futuros y agradecimientos"

$fin.Range("B17").Value = 43278
$fin.Range("B17").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C17").PasteSpecial(-4122)
$fin.Range("C17").Value = "añadidas las conclusiones y añadidas imágenes epipolares"

$fin.Range("B18").Value = 43278
$fin.Range("B18").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C18").PasteSpecial(-4122)
$fin.Range("C18").Value = "coreccion ortografica"

$fin.Range("B19").Value = 43279
$fin.Range("B19").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C19").PasteSpecial(-4122)
$fin.Range("C19").Value = "añadida explicacion sobre el output de libelas"

$fin.Range("B20").Value = 43280
$fin.Range("B20").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C20").PasteSpecial(-4122)
$fin.Range("C20").Value = "revision general y correcciones varias"

$fin.Range("B21").Value = 43280
$fin.Range("B21").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C21").PasteSpecial(-4122)
$fin.Range("C21").Value = "creacion del abstract"

$fin.Range("B22").Value = 43280
$fin.Range("B22").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C22").PasteSpecial(-4122)
$fin.Range("C22").Value = "añadidos resultados de roi"

$fin.Range("B23").Value = 43280
$fin.Range("B23").NumberFormat = "mm-dd-yy"
$ws1.Range("C5").Copy()
$fin.Range("C23").PasteSpecial(-4122)
$fin.Range("C23").Value = "coreccion ortografica"

# ---------------------------------------------------------------------------
# 6) Selections / active sheet (must be last so "Informe Final" stays active)
# ---------------------------------------------------------------------------
$ws1.Range("C43").Select()
$ws2.Range("C43").Select()
$seg.Range("B11").Select()
$fin.Range("B24").Select()
